$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.400.24"
$ws.Range("E2").Value = "  +1.29%  "

$ws.Range("D3").Value = "2.913.29"
$ws.Range("E3").Value = "  +3.97%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "353.92"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = "  -0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.26"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("E6").Value = "  +0.75%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.561"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = "  +0.79%  "

$ws.Range("E8").Value = "  +0.03%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.630"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("E9").Value = "  -0.10%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.11"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = "  -0.46%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0865"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = "  +2.88%  "

$ws.Range("E12").Value = "  +0.52%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.86"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = "  -0.78%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.79"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("E14").Value = "  +0.42%  "

$ws.Range("D15").Value = "3.370.37"
$ws.Range("E15").Value = "  +3.99%  "

$ws.Range("E16").Value = "  +5.43%  "

$ws.Range("D17").Value = "2.920.65"
$ws.Range("E17").Value = "  +4.15%  "

$ws.Range("D18").Value = "52.396.38"
$ws.Range("E18").Value = "  +1.41%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.67"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = "  +0.70%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.31"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = "  +3.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.28"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("E21").Value = "  +2.61%  "

$ws.Range("D22").Value = "0.0₃0981"
$ws.Range("E22").Value = "  +0.88%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "70.81"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.38"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = "  +0.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.80"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = "  +1.24%  "

$ws.Range("E26").Value = "  +7.50%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.86"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = "  +2.67%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.00"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("E28").Value = "  -0.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.66"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("E29").Value = "  +3.00%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.03"
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = "  -3.14%  "

$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.62"
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = "  +7.63%  "

$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.33"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("E32").Value = "  +12.63%  "

$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0993"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("E33").Value = "  +11.71%  "

$ws.Range("E34").Value = "  +0.78%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "53.44"
$ws.Range("D35").NumberFormat = "General"
$ws.Range("E35").Value = "  +2.21%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0453"
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = "  +1.73%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.38"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("E38").Value = "  +6.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "19.04"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("E39").Value = "  +0.72%  "

$ws.Range("E40").Value = "  +2.89%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.81"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("E41").Value = "  +12.60%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.118"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = "  +2.19%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.27"
$ws.Range("D43").NumberFormat = "General"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "120.62"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = "  -0.39%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.61"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("E45").Value = "  +5.51%  "

$ws.Range("E46").Value = "  -2.31%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.53"
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = "  +3.44%  "

$ws.Range("D48").Value = "2.195.88"
$ws.Range("E48").Value = "  +4.20%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.264"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("E49").Value = "  +20.71%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0347"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = "  +12.94%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.960"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = "  +0.73%  "
